$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Sprint 3 Product Backlog values (Priority column + reordered stories) ---

# Row 2: Priority 0.5 -> 0 (Create a Trip - unchanged text)
$ws.Range("A2").Value = 0

# Rows 3-5: unchanged (Add Waypoints / View overview / Add notes)

# Row 6: Priority 1 -> 0.5 (Remove Waypoints from my trip - text unchanged)
$ws.Range("A6").Value = 0.5

# Row 7: Priority 1 -> 0.5 (Add Transportation to my trip - text unchanged)
$ws.Range("A7").Value = 0.5

# Row 8: Priority 1 -> 0.5 (Remove Transportation from my trip - text unchanged)
$ws.Range("A8").Value = 0.5

# Row 9: Priority 2 -> 1 (Add Lodging to a trip - text unchanged)
$ws.Range("A9").Value = 1

# Row 10: Priority 2 -> 1 (Remove Lodging from a trip - text unchanged)
$ws.Range("A10").Value = 1

# Row 11: Priority 2 -> 1.5, content becomes "View details of an item in my Trip"
$ws.Range("A11").Value = 1.5
$ws.Range("C11").Value = "View details of an item in my Trip"
$ws.Range("D11").Value = "I can remember details about an item (Waypoint, Transportation, Lodging) in my Trip including any custom notes I have added to the item"

# Row 12: Priority stays 2, content becomes "Update Waypoints in my trip"
$ws.Range("A12").Value = 2
$ws.Range("C12").Value = "Update Waypoints in my trip"
$ws.Range("D12").Value = "I can update where I will be at certain times during my trip"

# Row 13: Priority stays 2, content becomes "Update Transportation in my trip"
$ws.Range("A13").Value = 2
$ws.Range("C13").Value = "Update Transportation in my trip"
$ws.Range("D13").Value = "I can update how I plan to travel to/from Waypoints"

# Row 14: Priority stays 2, content becomes "Update Lodging in a trip"
$ws.Range("A14").Value = 2
$ws.Range("C14").Value = "Update Lodging in a trip"
$ws.Range("D14").Value = "I can update details on where I will be staying during a specified time period of my trip"

# --- View / window cosmetic changes ---
$excel.ActiveWindow.Zoom = 160
$ws.Range("C11").Select()
